# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 60 (pushing the existing
# rows 60-158 down to 61-159), and populated with the new observation's
# data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 60; everything currently at row 60
# onward shifts down by one row.
$ws.Rows("60").Insert()

# Populate the newly inserted row 60 with the new record's values.
$ws.Range("A60").Value = 8
$ws.Range("B60").Value = "Terminal La Palmera de La Serena"
$ws.Range("C60").Value = "Coquimbo"
$ws.Range("D60").Value = 44799
$ws.Range("E60").Value = 4
$ws.Range("F60").Value = 100112001
$ws.Range("G60").Value = "Berenjena"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 540
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 11000
$ws.Range("M60").Value = 10500
$ws.Range("N60").Value = "`$/caja 40 unidades"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 262
$ws.Range("Q60").Value = 40
$ws.Range("R60").Value = "Hortaliza"
